$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.192.26'
$ws.Cells.Item(2, 5).Value = '  -0.15%  '
$ws.Cells.Item(3, 4).Value = '1.848.80'
$ws.Cells.Item(3, 5).Value = '  -0.55%  '
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '246.23'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.97%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6978'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -1.87%  '
$ws.Cells.Item(7, 5).Value = '  +0.05%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07718'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.46%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3061'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.27%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.50'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.51%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07823'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.10%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '93.05'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.99%  '
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.130'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.59%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.843.72'
$ws.Cells.Item(14, 5).Value = '  -0.85%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6860'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -0.24%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.644'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +2.16%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000008315'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -1.33%  '
$ws.Cells.Item(18, 4).Value = '29.189.41'
$ws.Cells.Item(18, 5).Value = '  -0.12%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '241.34'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -3.35%  '
$ws.Cells.Item(20, 4).Value = '2.087.52'
$ws.Cells.Item(20, 5).Value = '  -1.13%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.73'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.89%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9998'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.520'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.05%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.08%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1509'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.98%  '
$ws.Cells.Item(26, 5).Value = '  -0.95%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.834'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.32%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.28'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -1.32%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.546'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -1.00%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.230'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.40%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.172'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.25%  '
$ws.Cells.Item(32, 5).Value = '  -0.36%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05121'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.59%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7968'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +5.00%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.869'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.41%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.147'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -1.50%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.692'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.45%  '
$ws.Cells.Item(38, 4).Value = '1.316.01'
$ws.Cells.Item(38, 5).Value = '  +7.67%  '
$ws.Cells.Item(39, 5).Value = '  +0.66%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.715'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.51%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9433'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +5.18%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.017'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +8.06%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '107.15'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -2.46%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.07%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.737'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.23%  '
$ws.Cells.Item(46, 5).Value = '  -0.69%  '
$ws.Cells.Item(47, 4).Value = '1.988.71'
$ws.Cells.Item(47, 5).Value = '  -0.88%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5183'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.10%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.11'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.88%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.764'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.92%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.995'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.10%  '
